$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 to I1:J1, then set header labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I2:J20 data values (unstyled, matching H column formatting)
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 8
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 7
$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 8
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 8
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 8
$ws.Range("I11").Value = 8
$ws.Range("J11").Value = 8
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 5
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 4
$ws.Range("I14").Value = 7
$ws.Range("J14").Value = 7
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 3
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 8
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 6
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 8
$ws.Range("I19").Value = 8
$ws.Range("J19").Value = 8
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 6
